$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the description text in D5: drop the comma after "scope"
$ws.Range("D5").Value = "Subject completes response to perturbation having steered the vehicle back to the center of the lane. Normally this would be tagged with temporal scope but avoiding definitions here."

# Move the active selection from E4 to D5
$ws.Range("D5").Select()
